# Equip.xlsx: unify the conception of DataNode, DataTable, Entity.
# The sheet formerly named "Property1" is renamed to "DataNode" to match
# the unified naming convention across the data-config workbooks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the (only) worksheet from "Property1" to "DataNode".
$ws.Name = "DataNode"

# Re-apply the sheet's on-screen selection (frozen bottom-left pane was
# left focused on cell L39 at save time in the source commit).
$ws.Activate()
$ws.Range("L39").Select()
